$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 191; existing rows 191..289 shift down to 192..290.
$ws.Rows.Item(191).Insert()

# Populate the newly inserted row 191 with the new record.
$ws.Cells.Item(191, 1).Value = 10
$ws.Cells.Item(191, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(191, 3).Value = "La Araucanía"
$ws.Cells.Item(191, 4).Value = 45089
$ws.Cells.Item(191, 5).Value = 9
$ws.Cells.Item(191, 6).Value = "Fruta"
$ws.Cells.Item(191, 7).Value = 100104
$ws.Cells.Item(191, 8).Value = "Frutos de pepita"
$ws.Cells.Item(191, 9).Value = 100104003
$ws.Cells.Item(191, 10).Value = "Membrillo"
$ws.Cells.Item(191, 11).Value = "Champion"
$ws.Cells.Item(191, 12).Value = "Primera"
$ws.Cells.Item(191, 13).Value = 150
$ws.Cells.Item(191, 14).Value = 14000
$ws.Cells.Item(191, 15).Value = 14000
$ws.Cells.Item(191, 16).Value = 14000
$ws.Cells.Item(191, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(191, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(191, 19).Value = 778
$ws.Cells.Item(191, 20).Value = 18

# Match the date-formatted style used by column D elsewhere in the sheet.
$ws.Cells.Item(191, 4).NumberFormat = $ws.Cells.Item(192, 4).NumberFormat
